$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 171; existing rows 171:294 shift down to 172:295
$ws.Rows("171:171").Insert()

# Populate the newly inserted row 171 with the new weekly record
$ws.Range("A171").Value = 7
$ws.Range("B171").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C171").Value = "Ñuble"
$ws.Range("D171").Value = 44777
$ws.Range("E171").Value = 16
$ws.Range("F171").Value = 100112008
$ws.Range("G171").Value = "Coliflor"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 240
$ws.Range("K171").Value = 950
$ws.Range("L171").Value = 1000
$ws.Range("M171").Value = 975
$ws.Range("N171").Value = "`$/unidad"
$ws.Range("O171").Value = "Provincia de Diguillín"
$ws.Range("P171").Value = 975
$ws.Range("Q171").Value = 1
$ws.Range("R171").Value = "Hortaliza"
